$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-15 from 45207 to 45208
$ws.Range("C2:C15").Value = 45208
